$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.554.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.98"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5183"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3937"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07703"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.96"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.114"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.02"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.284"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.003"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.550"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.825.14"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.48"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001081"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06608"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.058"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.558.50"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.14"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.242"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.45%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.68"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.63"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.037.49"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.424"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.09"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.137"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.651"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.652"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07231"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2246"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.988"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02335"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.158"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6251"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.186"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.45"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5917"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.717"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.70"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.983"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.36%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06939"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.81%  "
